{"js": "// Word JS API (Office.js) script \u2014 body of `async (context) => { ... }`\n//\n// Target edit (from the authoritative OOXML diff):\n//   1. Title date run: \"2025-07-15 Tuesday\" -> \"2025-07-16 Wednesday\"\n//   2. The worksheet table has 20 rows (5 \"content\" rows holding the\n//      division problems, each followed by 3 blank rows). Every one of\n//      the 25 problem cells gets new text, in place (row/column count,\n//      run formatting (rFonts/sz) and paragraph formatting (jc) are\n//      unchanged \u2014 only the `w:t` text differs).\n//\n// Because several of the new cell values are identical to OLD values\n// that live elsewhere in the table (e.g. cell (0,0) becomes \"46\u00f76=\",\n// which is the pre-edit text of cell (4,4)), a blind document-wide\n// \"replace all occurrences of X\" is unsafe. Instead we address each\n// cell positionally with `table.getCell(row, col)` and then scope the\n// text search/replace to that single cell's body, which both avoids\n// any cross-cell collisions and preserves the existing run/paragraph\n// formatting (searching + replacing the found Range keeps its rPr/pPr,\n// unlike replacing the whole cell body).\n\n// 1) Title paragraph date text.\nconst titleResults = context.document.body.search(\"2025-07-15 Tuesday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"2025-07-16 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Table cell values, addressed by (row, col) to stay unambiguous.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> [newValue for col0..col4]\nconst newRowValues = {\n  0: [\"46\u00f76=\", \"45\u00f73=\", \"39\u00f77=\", \"15\u00f75=\", \"46\u00f74=\"],\n  4: [\"22\u00f78=\", \"87\u00f79=\", \"80\u00f72=\", \"81\u00f78=\", \"46\u00f78=\"],\n  8: [\"36\u00f76=\", \"74\u00f77=\", \"31\u00f79=\", \"78\u00f73=\", \"31\u00f75=\"],\n  12: [\"70\u00f74=\", \"90\u00f76=\", \"68\u00f79=\", \"37\u00f78=\", \"37\u00f75=\"],\n  16: [\"67\u00f72=\", \"70\u00f72=\", \"68\u00f72=\", \"44\u00f76=\", \"59\u00f76=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newValues = newRowValues[rowIndex];\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.body.load(\"text\");\n    await context.sync();\n\n    const oldText = cell.body.text.replace(/\\r/g, \"\").trim();\n    const newText = newValues[col];\n    if (oldText === newText) {\n      continue; // already correct, nothing to do\n    }\n\n    const cellResults = cell.body.search(oldText, { matchCase: true });\n    cellResults.load(\"items\");\n    await context.sync();\n\n    if (cellResults.items.length > 0) {\n      cellResults.items[0].insertText(newText, Word.InsertLocation.replace);\n    } else {\n      // Fallback: whole-body replace (loses formatting only if the\n      // targeted search could not locate the run, which should not\n      // happen given the source document).\n      cell.body.insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Target edit (from the authoritative OOXML diff):\n#   1. Title date run: \"2025-07-15 Tuesday\" -> \"2025-07-16 Wednesday\"\n#   2. The worksheet table has 20 rows (5 \"content\" rows holding the\n#      division problems, each followed by 3 blank rows). Every one of\n#      the 25 problem cells gets new text, in place (row/column count,\n#      run formatting (font/size) and paragraph alignment are\n#      unchanged - only the text differs).\n#\n# Several new cell values collide with OLD values living elsewhere in\n# the table (e.g. row 1 col 1 becomes \"46\u00f76=\", which is the pre-edit\n# text of row 5 col 5), so a single document-wide Find/Replace-all is\n# unsafe. Each cell is therefore addressed positionally via\n# Tables.Item / Rows.Item / Cells.Item, and updated by assigning\n# Cell.Range.Text directly - this reuses the cell's existing run, so\n# its font/size and the paragraph's alignment survive untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph date text - plain, unique text so Find/Replace is safe.\n$titleRange = $d.Content\n$titleRange.Find.Execute(\"2025-07-15 Tuesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-07-16 Wednesday\", 2) | Out-Null\n\n# 2) Table cell values, addressed by (row, col) to stay unambiguous.\n$table = $d.Tables.Item(1)\n\n$newRowValues = @{\n    1  = @(\"46\u00f76=\", \"45\u00f73=\", \"39\u00f77=\", \"15\u00f75=\", \"46\u00f74=\")\n    5  = @(\"22\u00f78=\", \"87\u00f79=\", \"80\u00f72=\", \"81\u00f78=\", \"46\u00f78=\")\n    9  = @(\"36\u00f76=\", \"74\u00f77=\", \"31\u00f79=\", \"78\u00f73=\", \"31\u00f75=\")\n    13 = @(\"70\u00f74=\", \"90\u00f76=\", \"68\u00f79=\", \"37\u00f78=\", \"37\u00f75=\")\n    17 = @(\"67\u00f72=\", \"70\u00f72=\", \"68\u00f72=\", \"44\u00f76=\", \"59\u00f76=\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $row = $table.Rows.Item($rowIndex)\n    $values = $newRowValues[$rowIndex]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $cell = $row.Cells.Item($col)\n        $newText = $values[$col - 1]\n        $cell.Range.Text = $newText\n    }\n}\n"}
